$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that can be safely set as text via Value2 without becoming numbers
# (non-numeric strings, or numeric-looking strings that Excel cannot parse as a number,
#  such as thousand-grouped "29.480.68")
$textUpdates = @{
    'D2' = '29.480.68'
    'E2' = '  +0.00%  '
    'D3' = '1.852.46'
    'E3' = '  +0.00%  '
    'E4' = '  -0.11%  '
    'E5' = '  -0.70%  '
    'E6' = '  +0.21%  '
    'E7' = '  -0.08%  '
    'E8' = '  +0.13%  '
    'E9' = '  -0.63%  '
    'E10' = '  +1.25%  '
    'E11' = '  -0.10%  '
    'D12' = '1.858.07'
    'E12' = '  -2.29%  '
    'E13' = '  -0.13%  '
    'E14' = '  -0.78%  '
    'E15' = '  +4.21%  '
    'E16' = '  -0.15%  '
    'B17' = 'Uniswap'
    'C17' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'E17' = '  -0.40%  '
    'B18' = 'WrappedBTC'
    'C18' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D18' = '29.509.54'
    'E18' = '  -0.32%  '
    'B19' = 'BitcoinCash'
    'C19' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'E19' = '  -1.95%  '
    'B20' = 'Avalanche'
    'C20' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'E20' = '  -0.29%  '
    'B21' = 'Dai'
    'C21' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'E21' = '  -0.14%  '
    'B22' = 'Chainlink'
    'C22' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E22' = '  -1.78%  '
    'B23' = 'BinanceUSD'
    'C23' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'E23' = '  -0.12%  '
    'B24' = 'Monero'
    'C24' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E24' = '  +0.91%  '
    'B25' = 'Stellar'
    'C25' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'E25' = '  +0.29%  '
    'B26' = 'Cosmos'
    'C26' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E26' = '  -1.16%  '
    'B27' = 'EthereumClassic'
    'C27' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'E27' = '  -0.35%  '
    'B28' = 'PancakeSwap'
    'C28' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'E28' = '  -0.61%  '
    'B29' = 'Toncoin'
    'C29' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'E29' = '  +1.27%  '
    'B30' = 'Hedera'
    'C30' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E30' = '  -1.93%  '
    'B31' = 'Filecoin'
    'C31' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E31' = '  +0.05%  '
    'B32' = 'InternetComputer(DFINITY)'
    'C32' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E32' = '  +0.36%  '
    'B33' = 'LidoDAOToken'
    'C33' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'E33' = '  -2.34%  '
    'B34' = 'ARBITRUM'
    'C34' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'E34' = '  -0.83%  '
    'B35' = 'ImmutableX'
    'C35' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'E35' = '  -0.97%  '
    'B36' = 'HuobiToken'
    'C36' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'E36' = '  +0.03%  '
    'B37' = 'Maker'
    'C37' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D37' = '1.250.32'
    'E37' = '  -0.45%  '
    'B38' = 'VeChain'
    'C38' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E38' = '  +0.26%  '
    'B39' = 'MXToken'
    'C39' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'E39' = '  -0.80%  '
    'B40' = 'TrustWalletToken'
    'C40' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E40' = '  +0.31%  '
    'B41' = 'FraxShare'
    'C41' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E41' = '  +0.92%  '
    'B42' = 'PaxDollar'
    'C42' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'E42' = '  -0.02%  '
    'B43' = 'Quant'
    'C43' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E43' = '  +0.18%  '
    'B44' = 'Aave'
    'C44' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'E44' = '  -2.05%  '
    'B45' = 'BabyDogeCoin'
    'C45' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'E45' = '  -4.62%  '
    'B46' = 'Aptos'
    'C46' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'E46' = '  -3.63%  '
    'B47' = 'EnergySwap'
    'C47' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E47' = '  -0.20%  '
    'B49' = 'RenderToken'
    'C49' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'E49' = '  -0.38%  '
    'B50' = 'Algorand'
    'C50' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E50' = '  +0.57%  '
    'B51' = 'Cronos'
    'C51' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'E51' = '  -0.64%  '
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value2 = $textUpdates[$ref]
}

# Cells whose new text would be auto-parsed as a number by Excel (e.g. "1.000", "0.9989").
# Force text format before assigning, then clear the style back to Normal so no residual
# number-format styling is left on the cell (matches original unstyled inline-string cells).
$guardedUpdates = @{
    'D4' = '0.9989'
    'D5' = '240.68'
    'D6' = '0.6324'
    'D7' = '1.000'
    'D9' = '0.2964'
    'D10' = '24.79'
    'D11' = '0.07706'
    'D13' = '5.005'
    'D14' = '0.6842'
    'D15' = '0.00001028'
    'D16' = '83.39'
    'D17' = '6.157'
    'D19' = '229.18'
    'D20' = '12.52'
    'D21' = '1.0000'
    'D22' = '7.551'
    'D23' = '1.000'
    'D24' = '156.79'
    'D25' = '0.1402'
    'D26' = '8.387'
    'D27' = '17.68'
    'D28' = '1.469'
    'D29' = '1.271'
    'D30' = '0.05677'
    'D31' = '4.133'
    'D32' = '4.037'
    'D33' = '1.849'
    'D34' = '1.161'
    'D35' = '0.7175'
    'D36' = '2.593'
    'D38' = '0.01810'
    'D39' = '2.779'
    'D40' = '0.9090'
    'D41' = '6.202'
    'D42' = '1.000'
    'D43' = '101.87'
    'D44' = '66.56'
    'D45' = '0.00000000120'
    'D46' = '7.107'
    'D47' = '9.152'
    'D48' = '0.4031'
    'D49' = '1.707'
    'D50' = '0.1125'
    'D51' = '0.05723'
}
foreach ($ref in $guardedUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value2 = $guardedUpdates[$ref]
    $cell.Style = "Normal"
}
